$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G9").Value = 41
$ws.Range("G10").Value = 41
$ws.Range("G11").Value = 41
$ws.Range("G12").Value = 41
